$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatasetTable")

# microclimUS is row 8: set Wind (K) and Snow (O) to "T"
$ws.Range("K8").Value = "T"
$ws.Range("O8").Value = "T"

$ws.Range("O8").Select()
